$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert the new 'College' column at C, shifting Email/TotalAttendance/Feedback right ---
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "College"

# --- Step 2: fill in College values for the existing rows 2-5 ---
$ws.Range("C2").Value = "ABC College"
$ws.Range("C3").Value = "XYZ Institute"
$ws.Range("C4").Value = "ABC College"
$ws.Range("C5").Value = "LMN University"

# --- Step 3: append new student rows 6-16 ---
$ws.Range("A6").Value = "Charlie Ray"
$ws.Range("B6").Value = "'105"
$ws.Range("C6").Value = "XYZ Institute"
$ws.Range("D6").Value = "charlie@example.com"
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = "Can do better"
$ws.Range("A7").Value = "Diana Prince"
$ws.Range("B7").Value = "'106"
$ws.Range("C7").Value = "ABC College"
$ws.Range("D7").Value = "diana@example.com"
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = "Outstanding"
$ws.Range("A8").Value = "Ethan Hunt"
$ws.Range("B8").Value = "'107"
$ws.Range("C8").Value = "LMN University"
$ws.Range("D8").Value = "ethan@example.com"
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = "Irregular"
$ws.Range("A9").Value = "Fiona Glen"
$ws.Range("B9").Value = "'108"
$ws.Range("C9").Value = "XYZ Institute"
$ws.Range("D9").Value = "fiona@example.com"
$ws.Range("E9").Value = 14
$ws.Range("F9").Value = "Consistent"
$ws.Range("A10").Value = "George White"
$ws.Range("B10").Value = "'109"
$ws.Range("C10").Value = "ABC College"
$ws.Range("D10").Value = "george@example.com"
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = "Good effort"
$ws.Range("A11").Value = "Hannah Moore"
$ws.Range("B11").Value = "'110"
$ws.Range("C11").Value = "LMN University"
$ws.Range("D11").Value = "hannah@example.com"
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = "Well done"
$ws.Range("A12").Value = "Ian Scott"
$ws.Range("B12").Value = "'111"
$ws.Range("C12").Value = "XYZ Institute"
$ws.Range("D12").Value = "ian@example.com"
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = "Punctual"
$ws.Range("A13").Value = "Julia Chen"
$ws.Range("B13").Value = "'112"
$ws.Range("C13").Value = "ABC College"
$ws.Range("D13").Value = "julia@example.com"
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = "Could improve"
$ws.Range("A14").Value = "Kevin Brooks"
$ws.Range("B14").Value = "'113"
$ws.Range("C14").Value = "LMN University"
$ws.Range("D14").Value = "kevin@example.com"
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "Needs attention"
$ws.Range("A15").Value = "Laura King"
$ws.Range("B15").Value = "'114"
$ws.Range("C15").Value = "XYZ Institute"
$ws.Range("D15").Value = "laura@example.com"
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = "Steady progress"
$ws.Range("A16").Value = "Michael Roy"
$ws.Range("B16").Value = "'115"
$ws.Range("C16").Value = "ABC College"
$ws.Range("D16").Value = "michael@example.com"
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "Very active"
